# Auto-generated edit script: refresh market-price derived values
# across the Asura_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1175
$ws.Range("I18").Value = 1175
$ws.Range("K18").Value = 1175
$ws.Range("M18").Value = -891
$ws.Range("H129").Value = 1079.4482
$ws.Range("J129").Value = 1122.537
$ws.Range("L129").Value = 3367.611
$ws.Range("N129").Value = -13367.611

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10378.4795
$ws.Range("I32").Value = 11207.8125
$ws.Range("J32").Value = 4481
$ws.Range("K32").Value = 11207.8125
$ws.Range("L32").Value = 4481
$ws.Range("M32").Value = -10920.8125
$ws.Range("N32").Value = -5055
$ws.Range("H123").Value = 24171
$ws.Range("J123").Value = 24171
$ws.Range("L123").Value = 24171
$ws.Range("N123").Value = -33971

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 31250
$ws.Range("J17").Value = 50000
$ws.Range("L17").Value = 50000
$ws.Range("N17").Value = -50348
$ws.Range("H19").Value = 197.14285
$ws.Range("I19").Value = 197.14285
$ws.Range("K19").Value = 197.14285
$ws.Range("M19").Value = -27.14285000000001
$ws.Range("H24").Value = 197.14285
$ws.Range("I24").Value = 197.14285
$ws.Range("K24").Value = 197.14285
$ws.Range("M24").Value = -27.14285000000001
$ws.Range("H31").Value = 1869.1364
$ws.Range("I31").Value = 1699.9412
$ws.Range("J31").Value = 2444.4
$ws.Range("K31").Value = 1699.9412
$ws.Range("L31").Value = 2444.4
$ws.Range("M31").Value = -1404.9412
$ws.Range("N31").Value = -3034.4
$ws.Range("H34").Value = 1869.1364
$ws.Range("I34").Value = 1699.9412
$ws.Range("J34").Value = 2444.4
$ws.Range("K34").Value = 1699.9412
$ws.Range("L34").Value = 2444.4
$ws.Range("M34").Value = -1497.9412
$ws.Range("N34").Value = -2848.4
$ws.Range("H41").Value = 1819.6666
$ws.Range("I41").Value = 1819.6666
$ws.Range("K41").Value = 1819.6666
$ws.Range("M41").Value = -1391.6666
$ws.Range("H42").Value = 50000
$ws.Range("I42").Value = 50000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 50000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -49407
$ws.Range("H132").Value = 399255.56
$ws.Range("I132").Value = 451796.38
$ws.Range("J132").Value = 5199.5
$ws.Range("K132").Value = 1355389.14
$ws.Range("L132").Value = 15598.5
$ws.Range("M132").Value = -1352859.14
$ws.Range("N132").Value = -20658.5
$ws.Range("H134").Value = 1605.8
$ws.Range("I134").Value = 1232.359
$ws.Range("K134").Value = 3697.077
$ws.Range("M134").Value = -1162.077
$ws.Range("H135").Value = 66835
$ws.Range("J135").Value = 66835
$ws.Range("L135").Value = 66835
$ws.Range("N135").Value = -76975
$ws.Range("H138").Value = 90780
$ws.Range("J138").Value = 90780
$ws.Range("L138").Value = 90780
$ws.Range("N138").Value = -101060
$ws.Range("H140").Value = 74488.57000000001
$ws.Range("J140").Value = 74488.57000000001
$ws.Range("L140").Value = 74488.57000000001
$ws.Range("N140").Value = -84848.57000000001
$ws.Range("N42").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3195.4546
$ws.Range("J109").Value = 3738.889
$ws.Range("L109").Value = 11216.667
$ws.Range("N109").Value = -13296.667
$ws.Range("H133").Value = 3660.4092
$ws.Range("J133").Value = 6967.25
$ws.Range("L133").Value = 20901.75
$ws.Range("N133").Value = -31021.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 22500
$ws.Range("J74").Value = 22500
$ws.Range("L74").Value = 22500
$ws.Range("N74").Value = -24372
$ws.Range("H77").Value = 22500
$ws.Range("J77").Value = 22500
$ws.Range("L77").Value = 67500
$ws.Range("N77").Value = -76860
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744
$ws.Range("H132").Value = 1577.6471
$ws.Range("I132").Value = 1230.4193
$ws.Range("J132").Value = 5165.6665
$ws.Range("K132").Value = 3691.2579
$ws.Range("L132").Value = 15496.9995
$ws.Range("M132").Value = -1161.2579
$ws.Range("N132").Value = -20556.9995
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H137").Value = 99000
$ws.Range("J137").Value = 99000
$ws.Range("L137").Value = 99000
$ws.Range("N137").Value = -109200
$ws.Range("H140").Value = 87090
$ws.Range("J140").Value = 87090
$ws.Range("L140").Value = 87090
$ws.Range("N140").Value = -97450
$ws.Range("N133").ClearContents()
$ws.Range("N135").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 1783.0333
$ws.Range("I136").Value = 1659.64
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 4978.92
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -2428.92
$ws.Range("N136").Value = -12300
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 82071.5
$ws.Range("J138").Value = 82071.5
$ws.Range("L138").Value = 82071.5
$ws.Range("N138").Value = -92351.5
$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -44860
$ws.Range("H141").Value = 42277.285
$ws.Range("J141").Value = 42277.285
$ws.Range("L141").Value = 42277.285
$ws.Range("N141").Value = -52637.285
$ws.Range("N135").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("N139").ClearContents()
